$wb = $excel.ActiveWorkbook

# OFF sheet: Week 16 logged + season sim updates for Home row
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 394
$wsOff.Range("C2").Value = 287
$wsOff.Range("D2").Value = 100
$wsOff.Range("E2").Value = 46

# DEF sheet: Week 16 logged + season sim updates for Home row
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 598
$wsDef.Range("C2").Value = 439
$wsDef.Range("D2").Value = 113
$wsDef.Range("E2").Value = 52
